$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates (Price and Volume(1h) columns, and row 26/27 coin swap)
$updates = [ordered]@{
    'D2' = '58.717.38'
    'E2' = '  -5.32%  '
    'D3' = '2.485.63'
    'E3' = '  -3.59%  '
    'D5' = '532.83'
    'E5' = '  -3.17%  '
    'D6' = '142.86'
    'E6' = '  -8.10%  '
    'E7' = '  -0.17%  '
    'D8' = '0.570'
    'E8' = '  -4.48%  '
    'D9' = '2.514.55'
    'E9' = '  -2.69%  '
    'D10' = '0.0998'
    'E10' = '  -4.53%  '
    'E11' = '  -2.73%  '
    'D12' = '5.54'
    'E12' = '  +0.83%  '
    'D13' = '0.350'
    'E13' = '  -4.71%  '
    'D14' = '2.933.76'
    'E14' = '  -3.31%  '
    'D15' = '23.73'
    'E15' = '  -7.92%  '
    'D16' = '58.667.86'
    'E16' = '  -5.20%  '
    'E17' = '  -5.23%  '
    'D18' = '2.507.31'
    'E18' = '  -2.88%  '
    'D19' = '11.25'
    'E19' = '  -3.32%  '
    'D20' = '4.26'
    'E20' = '  -6.85%  '
    'D21' = '321.07'
    'E21' = '  -5.13%  '
    'D22' = '0.996'
    'E22' = '  -0.37%  '
    'D23' = '5.72'
    'E23' = '  -5.47%  '
    'D24' = '60.64'
    'E24' = '  -4.64%  '
    'D25' = '0.437'
    'E25' = '  -11.72%  '
    'B26' = 'Kaspa'
    'C26' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D26' = '0.162'
    'E26' = '  -3.73%  '
    'B27' = 'WrappedeETH'
    'C27' = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
    'D27' = '2.614.91'
    'E27' = '  -3.03%  '
    'E28' = '  -0.36%  '
    'D29' = '7.70'
    'E29' = '  -6.05%  '
    'D30' = '6.82'
    'E30' = '  -7.51%  '
    'D31' = '0.0₃0770'
    'E31' = '  -8.64%  '
    'E32' = '  -8.23%  '
    'D33' = '1.77'
    'E33' = '  -6.76%  '
    'D34' = '0.997'
    'E34' = '  -0.21%  '
    'D35' = '155.79'
    'E35' = '  -4.15%  '
    'D36' = '1.41'
    'E36' = '  -2.55%  '
    'D37' = '18.45'
    'E37' = '  -4.19%  '
    'D38' = '4.36'
    'E38' = '  -10.76%  '
    'D39' = '1.59'
    'E39' = '  -11.70%  '
    'D40' = '5.83'
    'E40' = '  -3.69%  '
    'D41' = '305.85'
    'E41' = '  -7.83%  '
    'D42' = '36.73'
    'E42' = '  -2.23%  '
    'D43' = '3.65'
    'E43' = '  -7.73%  '
    'D44' = '0.799'
    'E44' = '  -13.03%  '
    'D46' = '0.593'
    'E46' = '  -2.30%  '
    'D47' = '10.77'
    'E47' = '  -1.44%  '
    'D48' = '124.39'
    'E48' = '  +0.89%  '
    'D49' = '0.0924'
    'E49' = '  -4.59%  '
    'D50' = '18.51'
    'E50' = '  -5.62%  '
    'D51' = '0.0513'
    'E51' = '  -6.73%  '
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    if ($cell.StartsWith('D')) {
        # Force text storage so numeric-looking strings (e.g. '532.83')
        # are not auto-converted to numbers by Excel, then restore the
        # cell's original (default) formatting/style.
        $range.NumberFormat = "@"
        $range.Value = $updates[$cell]
        $range.ClearFormats()
    } else {
        $range.Value = $updates[$cell]
    }
}

